# Read Data from excel in Login Validation
# Rename the default sheet to "Login" and populate it with test data
# used by the login-validation selenium tests.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Login"

# Header row
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"

# Test data rows (mix of invalid credentials plus the valid demo account)
$ws.Range("B3").Value = "12sdd35142"
$ws.Range("B4").Value = "1sdsd5142"
$ws.Range("A4").Value = "sdsddcdsad"
$ws.Range("A3").Value = "sdgdddssad"
$ws.Range("A2").Value = "sdddsdsad"
$ws.Range("A5").Value = "tomsmith"
$ws.Range("B5").Value = "SuperSecretPassword!"
$ws.Range("B2").Value = 1241235142

# Column widths
$ws.Columns.Item(1).ColumnWidth = 13.28515625
$ws.Columns.Item(2).ColumnWidth = 24.5703125

# Header formatting: bold-ish larger font, left aligned with indent
$headerRange = $ws.Range("A1:B1")
$headerRange.Font.Name = "Segoe UI"
$headerRange.Font.Size = 12
$headerRange.HorizontalAlignment = -4131
$headerRange.IndentLevel = 1

# Data formatting: left aligned with indent
$dataRange = $ws.Range("A2:B5")
$dataRange.HorizontalAlignment = -4131
$dataRange.IndentLevel = 1

$ws.PageSetup.Orientation = 1

$ws.Range("L9").Select()
